$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 6
$ws.Range("F3").Value = -1
$ws.Range("F9").Value = 2
$ws.Range("F10").Value = -2
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = 1
$ws.Range("F21").Value = -2
$ws.Range("F23").Value = -3
$ws.Range("F25").Value = 3
$ws.Range("F30").Value = -4
$ws.Range("F33").Value = -2
$ws.Range("F34").Value = -1
$ws.Range("F38").Value = -3
$ws.Range("F43").Value = -1
$ws.Range("F44").Value = 0
$ws.Range("F46").Value = 2
$ws.Range("F48").Value = 1
$ws.Range("F51").Value = -2
$ws.Range("F52").Value = 4
$ws.Range("F54").Value = 2
$ws.Range("F55").Value = 3
$ws.Range("F57").Value = 2
$ws.Range("F74").Value = 3
$ws.Range("F78").Value = -7
$ws.Range("F79").Value = -3
$ws.Range("F80").Value = 1
$ws.Range("F85").Value = -8
